# Update the "nome"/"cognome" value from "a" to "b" and the "email" value
# from "a@a.it" to "b@b.it" (post DB update correction), as described by
# the commit message "risoluzione aggiornamento dati post UPDATE DB".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the single data record:
#   B2 = nome, C2 = cognome, D2 = email
$ws.Range("B2").Value = "b"
$ws.Range("C2").Value = "b"
$ws.Range("D2").Value = "b@b.it"
